# Work Profile and new tenant support
# Adds new sprint-history rows to the AMSIN, BETA and AMS sheets, and
# fixes up the previously-unstyled AMS row 36.

$wb = $excel.ActiveWorkbook

$TimeFmt = "YYYY-MM-DD HH:MM:SS"

function Add-SprintRow {
    param(
        $ws,
        [int]$Row,
        [string]$RunDate,
        [double]$RunTime,
        [string]$SprintName,
        [double]$Total,
        [double]$Pass,
        [double]$Fail,
        [double]$TimeTaken
    )

    # Column A holds a literal date-shaped string (not a real date) - the
    # leading apostrophe stops Excel's autodetection turning it into a
    # date serial number, same as the source data (t="inlineStr"). Reset
    # the style back to Normal afterwards so the quote-prefix marker left
    # by the apostrophe entry doesn't linger on the cell format.
    $ws.Cells.Item($Row, 1).ClearContents()
    $ws.Cells.Item($Row, 1).Value = "'" + $RunDate
    $ws.Cells.Item($Row, 1).Style = "Normal"

    # Column B is a real number formatted as a date-time string.
    $ws.Cells.Item($Row, 2).ClearContents()
    $ws.Cells.Item($Row, 2).NumberFormat = $TimeFmt
    $ws.Cells.Item($Row, 2).Value = $RunTime

    $ws.Cells.Item($Row, 3).ClearContents()
    $ws.Cells.Item($Row, 3).Value = $SprintName

    $ws.Cells.Item($Row, 4).ClearContents()
    $ws.Cells.Item($Row, 4).Value = $Total

    $ws.Cells.Item($Row, 5).ClearContents()
    $ws.Cells.Item($Row, 5).Value = $Pass

    $ws.Cells.Item($Row, 6).ClearContents()
    $ws.Cells.Item($Row, 6).Value = $Fail

    $ws.Cells.Item($Row, 7).ClearContents()
    $ws.Cells.Item($Row, 7).Value = $TimeTaken
}

# ---------------------------------------------------------------------
# AMSIN sheet: new rows 62-67
# ---------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-SprintRow $wsAmsin 62 "2023-03-09" 44994.61513378472   "174fstcycle"    124 124 0  1.9
Add-SprintRow $wsAmsin 63 "2023-03-10" 44995.79170776621   "174ffiinnalrun" 124 121 3  2.82
Add-SprintRow $wsAmsin 64 "2023-03-13" 44998.46589782408   "174finalrun"    124 122 2  1.92
Add-SprintRow $wsAmsin 65 "2023-03-30" 45015.71254118055   "175scndcyc"     124 120 4  2.61
Add-SprintRow $wsAmsin 66 "2023-03-31" 45016.44677490741   "175fnlrun"      124 119 5  2.03
Add-SprintRow $wsAmsin 67 "2023-04-12" 45028.5508363341    "176fstrtail"    124 105 19 5.66

# ---------------------------------------------------------------------
# BETA sheet: new rows 33-34
# ---------------------------------------------------------------------
$wsBeta = $wb.Worksheets.Item("BETA")

Add-SprintRow $wsBeta 33 "2023-03-13" 44998.53276665509 "174beta" 124 124 0 1.62
Add-SprintRow $wsBeta 34 "2023-03-31" 45016.5270653125  "175beta" 124 122 2 1.85

# ---------------------------------------------------------------------
# AMS sheet: row 36 refresh + new rows 37-38
# ---------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

Add-SprintRow $wsAms 36 "2023-02-20" 44977.83446886574 "live173" 124 123 1 1.78
Add-SprintRow $wsAms 37 "2023-03-13" 44998.83526253473 "174live" 124 124 0 1.57
Add-SprintRow $wsAms 38 "2023-03-31" 45016.79327719907 "175live" 124 124 0 1.62
